$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new rows ------------------------------------------------------

$ws.Range("A43").Value = 1041361
$ws.Range("B43").Value = "KENT PANGESTU"
$ws.Range("C43").Value = "TJKT"

$ws.Range("A44").Value = 1041362
$ws.Range("B44").Value = "TRI FANDI JUNIOR"
$ws.Range("C44").Value = "TJKT"

$ws.Range("A45").Value = 1041363
$ws.Range("B45").Value = "EKA DHANI"
$ws.Range("C45").Value = "TJKT"

# --- Update existing rows -----------------------------------------------

# Row 28: USUP RAHARJO, A.MD -> USUP RAHARJO, S.KOM
$ws.Range("B28").Value = "USUP RAHARJO, S.KOM"

# Row 30: RAHMADIKA SURYA S.KOM -> RAHMADIKA SURYA SETIAWAN, S.KOM
$ws.Range("B30").Value = "RAHMADIKA SURYA SETIAWAN, S.KOM"

# Row 41: MARCELLINO RADITIO -> MARCELLINO RADITIO, S.KOOM
$ws.Range("B41").Value = "MARCELLINO RADITIO, S.KOOM"

# Row 42: ANDIKA's type DKV -> TJKT
$ws.Range("C42").Value = "TJKT"

# --- Formatting / filter cleanup -------------------------------------------

# Remove the number-format style that was applied to column A (ID Guru)
$ws.Range("A2:A45").ClearFormats()

# Remove the autofilter entirely
$ws.AutoFilterMode = $false

# --- Selection ---------------------------------------------------------

$ws.Range("C22").Select()
